$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 3 to reflect the latest FlashScore data
$ws.Range("G3").Value = 1.45
$ws.Range("I3").Value = 6.5
$ws.Range("J3").Value = 2
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Z3").Value = 9.5
$ws.Range("AJ3").Value = 21
$ws.Range("AO3").Value = 7

# Remove the Al Fateh - Al Riyadh match row (row 8), shrinking the used range to A1:BD7
$ws.Rows("8:8").Delete()
